# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计" (so it becomes the
#    2nd tab, pushing all existing quarter sheets one slot to the right).
# 2. Populate the new sheet with the 2022-Q3 fund-holding table (same shape
#    as the other quarter sheets).
# 3. Update the "总计" (totals) sheet: add a new row for 2022-Q3 at the top
#    of the data (row 2) and push the rest of the rows down, which also
#    reveals a previously-missing 2020-Q4 row at the bottom (row 9).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create & position the new "2022-Q3" worksheet
# ---------------------------------------------------------------------
# NOTE: worksheet object references returned by Item(...) track a sheet's
# *position*, not a stable identity - once we insert a new sheet, any
# reference fetched beforehand can silently start pointing at the newly
# inserted sheet instead. So insert/rename first, and only resolve the
# sheet references we need (by name) afterwards.
$beforeSheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q3"

$totalsSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Copy the header-row (bold/bordered) style and the column-A (bold/bordered)
# style from the existing "2022-Q2" sheet so the new sheet matches the look
# of its siblings without inventing new style entries.
$q2Sheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$q2Sheet.Range("A2").Copy($newSheet.Range("A2:A7"))

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$fundRows = @(
    @(0, "159745", "国泰中证全指建筑材料ETF",     "7.92", "99.14", "2.85", "0.2257", 10),
    @(1, "004856", "广发中证全指建筑材料指数A",    "7.66", "93.74", "2.69", "0.2061", 10),
    @(2, "004857", "广发中证全指建筑材料指数C",    "6.12", "93.74", "2.69", "0.1646", 10),
    @(3, "516750", "富国中证全指建筑材料ETF",      "0.82", "98.46", "2.83", "0.0232", 10),
    @(4, "014344", "鹏华中证500指数增强A",         "1.14", "92.67", "1.67", "0.0190", 7),
    @(5, "014345", "鹏华中证500指数增强C",         "0.73", "92.67", "1.67", "0.0122", 7)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $fundRow = $fundRows[$i]

    $newSheet.Cells.Item($r, 1).Value = $fundRow[0]

    # Columns B, D, E, F, G hold numeric-looking text in the source data
    # (t="inlineStr"), so force text entry with a leading apostrophe, then
    # strip the resulting quote-prefix style so the cell keeps the default
    # (unstyled) look like its neighbours.
    $newSheet.Cells.Item($r, 2).Value = "'" + $fundRow[1]
    $newSheet.Cells.Item($r, 2).ClearFormats()

    $newSheet.Cells.Item($r, 3).Value = $fundRow[2]

    $newSheet.Cells.Item($r, 4).Value = "'" + $fundRow[3]
    $newSheet.Cells.Item($r, 4).ClearFormats()

    $newSheet.Cells.Item($r, 5).Value = "'" + $fundRow[4]
    $newSheet.Cells.Item($r, 5).ClearFormats()

    $newSheet.Cells.Item($r, 6).Value = "'" + $fundRow[5]
    $newSheet.Cells.Item($r, 6).ClearFormats()

    $newSheet.Cells.Item($r, 7).Value = "'" + $fundRow[6]
    $newSheet.Cells.Item($r, 7).ClearFormats()

    $newSheet.Cells.Item($r, 8).Value = $fundRow[7]
}

# ---------------------------------------------------------------------
# Step 2: rewrite the "总计" table (A1:D9) with 2022-Q3 added at the top
#         and 2020-Q4 revealed as the new last row
# ---------------------------------------------------------------------

# Give the new last row (A9) the same bold/bordered style as column A of
# the row above it, before (re)writing all the values.
$totalsSheet.Range("A8").Copy($totalsSheet.Range("A9"))

# index, 日期, 持有数量(只), 持有市值(亿元)
$totalsRows = @(
    @(0, "2022-Q3", 6, 0.65),
    @(1, "2022-Q2", 10, 0.79),
    @(2, "2022-Q1", 16, 5.09),
    @(3, "2021-Q4", 14, 2.8),
    @(4, "2021-Q3", 10, 2.33),
    @(5, "2021-Q2", 11, 2.18),
    @(6, "2021-Q1", 12, 2.82),
    @(7, "2020-Q4", 4, 0.89)
)

for ($i = 0; $i -lt $totalsRows.Length; $i++) {
    $r = $i + 2
    $totalsRow = $totalsRows[$i]
    $totalsSheet.Cells.Item($r, 1).Value = $totalsRow[0]
    $totalsSheet.Cells.Item($r, 2).Value = $totalsRow[1]
    $totalsSheet.Cells.Item($r, 3).Value = $totalsRow[2]
    $totalsSheet.Cells.Item($r, 4).Value = $totalsRow[3]
}
